$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.405.23'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '2.301.21'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.02'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.04'
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.20'
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.34'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.968'
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.31'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '2.649.48'
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").Value = '2.304.04'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '42.508.84'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.26'
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '276.17'
$ws.Range("E23").Value = '  +6.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  +20.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.78'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.86'
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.86'
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0874'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("E34").Value = '  +5.26%  '
$ws.Range("E35").Value = '  -10.42%  '
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0373'
$ws.Range("E37").Value = '  +5.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.60'
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("E39").Value = '  +3.18%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.17'
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '95.29'
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '82.73'
$ws.Range("E46").Value = '  +10.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.07'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.00'
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").Value = '1.591.50'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.05'
$ws.Range("E51").Value = '  -5.81%  '
